{"js": "// Update the two-digit x two-digit multiplication answers in the table\n// from the \"before\" values to the newly generated ones. Each answer\n// lives alone in its own paragraph/run, so a simple whole-document\n// search & replace (by exact old text) is unambiguous.\nconst replacements = [\n  [\"64\u00d741=2624\", \"90\u00d754=4860\"],\n  [\"60\u00d742=2520\", \"23\u00d762=1426\"],\n  [\"12\u00d776=912\", \"83\u00d718=1494\"],\n  [\"18\u00d725=450\", \"48\u00d772=3456\"],\n  [\"84\u00d798=8232\", \"66\u00d730=1980\"],\n  [\"16\u00d754=864\", \"92\u00d734=3128\"],\n  [\"59\u00d737=2183\", \"75\u00d721=1575\"],\n  [\"39\u00d731=1209\", \"70\u00d793=6510\"],\n  [\"16\u00d770=1120\", \"31\u00d763=1953\"],\n  [\"63\u00d761=3843\", \"93\u00d737=3441\"],\n  [\"75\u00d759=4425\", \"96\u00d723=2208\"],\n  [\"46\u00d768=3128\", \"97\u00d748=4656\"],\n  [\"44\u00d787=3828\", \"88\u00d799=8712\"],\n  [\"61\u00d760=3660\", \"70\u00d797=6790\"],\n  [\"77\u00d786=6622\", \"17\u00d780=1360\"],\n  [\"95\u00d714=1330\", \"84\u00d767=5628\"],\n  [\"35\u00d731=1085\", \"26\u00d719=494\"],\n  [\"84\u00d737=3108\", \"31\u00d797=3007\"],\n  [\"46\u00d714=644\", \"29\u00d727=783\"],\n  [\"71\u00d730=2130\", \"83\u00d737=3071\"],\n  [\"99\u00d770=6930\", \"60\u00d711=660\"],\n  [\"48\u00d776=3648\", \"47\u00d793=4371\"],\n  [\"61\u00d791=5551\", \"42\u00d787=3654\"],\n  [\"68\u00d719=1292\", \"13\u00d797=1261\"],\n  [\"21\u00d719=399\", \"68\u00d732=2176\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit x two-digit multiplication answers in the table\n# from the \"before\" values to the newly generated ones. Each answer\n# lives alone in its own paragraph/run, so a whole-document Find &\n# Replace (exact text match) for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"64\u00d741=2624\"; new=\"90\u00d754=4860\"},\n    @{old=\"60\u00d742=2520\"; new=\"23\u00d762=1426\"},\n    @{old=\"12\u00d776=912\";  new=\"83\u00d718=1494\"},\n    @{old=\"18\u00d725=450\";  new=\"48\u00d772=3456\"},\n    @{old=\"84\u00d798=8232\"; new=\"66\u00d730=1980\"},\n    @{old=\"16\u00d754=864\";  new=\"92\u00d734=3128\"},\n    @{old=\"59\u00d737=2183\"; new=\"75\u00d721=1575\"},\n    @{old=\"39\u00d731=1209\"; new=\"70\u00d793=6510\"},\n    @{old=\"16\u00d770=1120\"; new=\"31\u00d763=1953\"},\n    @{old=\"63\u00d761=3843\"; new=\"93\u00d737=3441\"},\n    @{old=\"75\u00d759=4425\"; new=\"96\u00d723=2208\"},\n    @{old=\"46\u00d768=3128\"; new=\"97\u00d748=4656\"},\n    @{old=\"44\u00d787=3828\"; new=\"88\u00d799=8712\"},\n    @{old=\"61\u00d760=3660\"; new=\"70\u00d797=6790\"},\n    @{old=\"77\u00d786=6622\"; new=\"17\u00d780=1360\"},\n    @{old=\"95\u00d714=1330\"; new=\"84\u00d767=5628\"},\n    @{old=\"35\u00d731=1085\"; new=\"26\u00d719=494\"},\n    @{old=\"84\u00d737=3108\"; new=\"31\u00d797=3007\"},\n    @{old=\"46\u00d714=644\";  new=\"29\u00d727=783\"},\n    @{old=\"71\u00d730=2130\"; new=\"83\u00d737=3071\"},\n    @{old=\"99\u00d770=6930\"; new=\"60\u00d711=660\"},\n    @{old=\"48\u00d776=3648\"; new=\"47\u00d793=4371\"},\n    @{old=\"61\u00d791=5551\"; new=\"42\u00d787=3654\"},\n    @{old=\"68\u00d719=1292\"; new=\"13\u00d797=1261\"},\n    @{old=\"21\u00d719=399\";  new=\"68\u00d732=2176\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $p.new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}\n"}
